$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8837209939956665
$ws.Range("B1").Value = 2.69945240020752
$ws.Range("C1").Value = 4.633284091949463
$ws.Range("D1").Value = 2.200158834457397
$ws.Range("E1").Value = 1.298065066337585
